$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (O2, P2, S2, T2)
$ws.Range("O2").Value = 0.7927501397588634
$ws.Range("P2").Value = 0.7927501397588635
$ws.Range("S2").Value = 0.7927501397588634
$ws.Range("T2").Value = 0.7927501397588635

# Update existing row 3 values (M3..T3)
$ws.Range("M3").Value = 0.8869683333333332
$ws.Range("N3").Value = 2.660905
$ws.Range("O3").Value = 0.2003151406163121
$ws.Range("P3").Value = 0.2003151406163121
$ws.Range("Q3").Value = 0.06499999602777777
$ws.Range("R3").Value = 0.5849999642499999
$ws.Range("S3").Value = 0.2003151406163121
$ws.Range("T3").Value = 0.2003151406163121

# Add new row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bmp7"
$ws.Range("C4").Value = "Bmpr1b"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07328333333333333
$ws.Range("H4").Value = 0.21985
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.030706
$ws.Range("N4").Value = 0.09211800000000001
$ws.Range("O4").Value = 0.006934719624824425
$ws.Range("P4").Value = 0.006934719624824427
$ws.Range("Q4").Value = 0.002250238033333333
$ws.Range("R4").Value = 0.0202521423
$ws.Range("S4").Value = 0.006934719624824425
$ws.Range("T4").Value = 0.006934719624824427
